# Update loading_percent values per case "Case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 8.919008135648394
    "C2" = 5.224702408197171
    "E2" = 12.65666093227336
    "F2" = 16.86991607391245
    "G2" = 3.637179407437194
    "K2" = 8.203292644561094
    "M2" = 13.52572405511802
    "N2" = 18.7286439175631
    "O2" = 22.10656030117271
    "B3" = 8.662430786648367
    "C3" = 5.123336841014574
    "E3" = 12.43270288794568
    "F3" = 15.89584955866815
    "G3" = 3.638829452682554
    "K3" = 8.035691987721417
    "M3" = 13.36373318392015
    "N3" = 18.79073583527846
    "O3" = 22.18856529323418
    "B4" = 8.502631965821045
    "C4" = 5.059441136266366
    "E4" = 12.29750614191613
    "F4" = 15.26997757108489
    "G4" = 3.639895973291792
    "K4" = 7.932351903172238
    "M4" = 13.26656067721565
    "N4" = 18.83064378844361
    "O4" = 22.24358680610488
    "B5" = 8.437049022182631
    "C5" = 5.033006851033083
    "E5" = 12.24307613344562
    "F5" = 15.008197319934
    "G5" = 3.640344055342174
    "K5" = 7.890189716716403
    "M5" = 13.22758391562424
    "N5" = 18.84735637855539
    "O5" = 22.26718046601945
    "B6" = 8.42613420084639
    "C6" = 5.028594138306354
    "E6" = 12.23408053558174
    "F6" = 14.96433081551589
    "G6" = 3.640419273584785
    "K6" = 7.883187421565226
    "M6" = 13.22115072596437
    "N6" = 18.85015869670284
    "O6" = 22.27116888745825
    "B7" = 8.501749225752816
    "C7" = 5.059086211113027
    "E7" = 12.29676928475996
    "F7" = 15.26647399323133
    "G7" = 3.639901961703028
    "K7" = 7.931783414267225
    "M7" = 13.26603244559262
    "N7" = 18.83086735736002
    "O7" = 22.24390025648781
    "B8" = 8.831069856563477
    "C8" = 5.190107345798806
    "E8" = 12.57900910086885
    "F8" = 16.5399640634477
    "G8" = 3.637737287975349
    "K8" = 8.145628333887059
    "M8" = 13.46942206412236
    "N8" = 18.74968391570397
    "O8" = 22.13386491268468
    "B9" = 9.454687985511194
    "C9" = 5.433050621074254
    "E9" = 13.14723557912611
    "F9" = 19.00274580682531
    "G9" = 3.633914050796732
    "K9" = 8.559063251261676
    "M9" = 13.8843775200867
    "N9" = 18.60456982986504
    "O9" = 21.95524702833527
    "B10" = 9.894143288057071
    "C10" = 5.601960357270194
    "E10" = 13.56891731672967
    "F10" = 20.67494806633232
    "G10" = 3.63135950318925
    "K10" = 8.856037907764062
    "M10" = 14.19625453519216
    "N10" = 18.50645232817607
    "O10" = 21.84681387040371
    "B11" = 10.08904477752217
    "C11" = 5.676520131961047
    "E11" = 13.76068159219125
    "F11" = 21.3917225636224
    "G11" = 3.630252039567652
    "K11" = 8.989047972031134
    "M11" = 14.3390570718402
    "N11" = 18.46364222390167
    "O11" = 21.80246501609464
    "B12" = 10.16205882279117
    "C12" = 5.704410976853926
    "E12" = 13.8332105580082
    "F12" = 21.65686569030329
    "G12" = 3.629840482024248
    "K12" = 9.03906777539412
    "M12" = 14.39321640345716
    "N12" = 18.44769203059638
    "O12" = 21.78638951988322
    "B13" = 10.14637016946272
    "C13" = 5.698419692650676
    "E13" = 13.81759528806346
    "F13" = 21.60004134736742
    "G13" = 3.629928771343436
    "K13" = 9.028311340511845
    "M13" = 14.38154931886026
    "N13" = 18.4511155986902
    "O13" = 21.78981967094903
    "B14" = 10.09506791821031
    "C14" = 5.678821695495684
    "E14" = 13.76665079977047
    "F14" = 21.4136618050453
    "G14" = 3.630218024077697
    "K14" = 8.993170376379362
    "M14" = 14.34351141247288
    "N14" = 18.46232476679565
    "O14" = 21.80112806161693
    "B15" = 10.06353882869422
    "C15" = 5.66677219725958
    "E15" = 13.73543207311252
    "F15" = 21.29868154950795
    "G15" = 3.630396216325892
    "K15" = 8.97159878925792
    "M15" = 14.32022143938384
    "N15" = 18.46922466783462
    "O15" = 21.80814841828189
    "B16" = 9.881299027369105
    "C16" = 5.597040489146337
    "E16" = 13.55637699722755
    "F16" = 20.62722412089977
    "G16" = 3.631432974205866
    "K16" = 8.847299222901942
    "M16" = 14.18693640783008
    "N16" = 18.50928666805958
    "O16" = 21.8498125673005
    "B17" = 9.768165424682232
    "C17" = 5.553667321900947
    "E17" = 13.44645952226042
    "F17" = 20.20408069597325
    "G17" = 3.632082951855718
    "K17" = 8.770476011119035
    "M17" = 14.1053720116806
    "N17" = 18.53432970875047
    "O17" = 21.87664903120265
    "B18" = 9.702626275749848
    "C18" = 5.52850666725953
    "E18" = 13.38323672078672
    "F18" = 19.95656407809801
    "G18" = 3.63246194450794
    "K18" = 8.726095764708869
    "M18" = 14.05854866110504
    "N18" = 18.54890556343169
    "O18" = 21.89255304796465
    "B19" = 9.680357786048614
    "C19" = 5.519951520882386
    "E19" = 13.36183283166784
    "F19" = 19.87204792380568
    "G19" = 3.632591149345833
    "K19" = 8.711037654492886
    "M19" = 14.04271214046162
    "N19" = 18.55387023502392
    "O19" = 21.8980182405789
    "B20" = 9.78025762802749
    "C20" = 5.558306689305057
    "E20" = 13.4581611313502
    "F20" = 20.24955283636154
    "G20" = 3.632013228653991
    "K20" = 8.778674366719546
    "M20" = 14.11404569901064
    "N20" = 18.53164606673474
    "O20" = 21.8737437476797
    "B21" = 10.11015862699388
    "C21" = 5.68458753997971
    "E21" = 13.78161743448706
    "F21" = 21.46857628470577
    "G21" = 3.630132851747975
    "K21" = 9.003501950267047
    "M21" = 14.35468220923308
    "N21" = 18.45902528695611
    "O21" = 21.79778699387099
    "B22" = 10.32112996573657
    "C22" = 5.765111695022359
    "E22" = 13.992467528221
    "F22" = 22.22866616901552
    "G22" = 3.628949450340969
    "K22" = 9.148389814838625
    "M22" = 14.51241379392213
    "N22" = 18.41308452739844
    "O22" = 21.75233373506873
    "B23" = 10.20897653780795
    "C23" = 5.722323056295757
    "E23" = 13.88000825300517
    "F23" = 21.82633154458858
    "G23" = 3.629576900214835
    "K23" = 9.071263190231015
    "M23" = 14.42820364025024
    "N23" = 18.43746519669569
    "O23" = 21.77620885811531
    "B24" = 9.774792285907015
    "C24" = 5.556209930457812
    "E24" = 13.45287092057952
    "F24" = 20.22900810905287
    "G24" = 3.632044733964353
    "K24" = 8.77496855161859
    "M24" = 14.1101241058521
    "N24" = 18.53285878585573
    "O24" = 21.87505574576806
    "B25" = 9.288921773720878
    "C25" = 5.368938477639198
    "E25" = 12.99243535450638
    "F25" = 18.34778573295695
    "G25" = 3.634903472394786
    "K25" = 8.448188267806792
    "M25" = 13.77068557623288
    "N25" = 18.64232829514189
    "O25" = 21.80814841828189
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
